$wb = $excel.ActiveWorkbook

# --- 1) Add new worksheet "Eligible  2s" as the last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsNew.Name = "Eligible  2s"

$wsNew.Range("B1").Value = "Eligible Two Year Olds"
$wsNew.Range("A2").Value = "Uptake measure"
$wsNew.Range("B2").Value = "add text here"
$wsNew.Range("A3").Value = "Demand measure"
$wsNew.Range("B3").Value = "add text here"
$wsNew.Range("A4").Value = "Relationship"
$wsNew.Range("B4").Value = "add text here"
$wsNew.Range("A5").Value = "Data sources"
$wsNew.Range("B5").Value = "add text here"
$wsNew.Range("A6").Value = "Publication year"
$wsNew.Range("B6").Value = "add text here"
$wsNew.Range("A7").Value = "Display year"
$wsNew.Range("B7").Value = "add text here"
$wsNew.Range("A8").Value = "Additional information"
$wsNew.Range("B8").Value = "add text here"

$wsNew.Columns.Item(1).ColumnWidth = 17.9140625
$wsNew.Columns.Item(2).ColumnWidth = 16.9140625

# --- 2) Update "expected" -> "estimated" wording across sheets ---

# General sheet: summary text, "Modelling expected uptake" heading, and methodology text
$wsGeneral = $wb.Worksheets.Item("General")
$wsGeneral.Range("B2").Value = "This dashboard aims to support users to identify local areas where there may be unmet need in relation to child poverty. Unmet need is defined as cases where families eligible for support, such as social security, do not access this. `n`nThe dashboard  highlights data zones that may have unmet need. That is, it highlights areas where the actual rates of uptake for benefits and support related to child poverty differ significantly from modelled estimated  uptake of these."
$wsGeneral.Range("A4").Value = "Modelling estimated uptake"
$wsGeneral.Range("B4").Value = "To model estimated uptake, pairs of measures are used, with each pairing including an “uptake” measure and a “demand” measure. The uptake measure is the benefit or service of interest, and the demand measure is a suitable comparator which can be used to estimate uptake. For example, uptake of Free School Meals for children aged 10 to 18 years is used as an uptake measure and the proportion of children aged 10 to 18 years in low-income families is used as a comparative measure of demand. Each pairing of uptake and demand measures have similar eligibility requirements and therefore a linear relationship, meaning that higher rates of demand are typically associated with higher rates of uptake. This relationship can be used to estimate demand.`n`n A linear regression model is fitted between each uptake measure and its associated demand measure. The regression model calculates for each data zone what the estimated uptake rate would be given the demand rate. Residual values - the difference between the observed rate and the estimated rate - are calculated and then standardised. Data zones where the standardised residual is above 2 or below minus 2 are then highlighted in the dashboard as areas where uptake differs from what is estimated.  `n`nWhen interpreting higher/lower than estimated uptake, the estimated value does not mean that this is what the rate should be but is the rate most likely at the associated demand rate within the model. In practice, uptake rates will be higher or lower than the modelled values for a range of reasons. "

# Universal Credit sheet: additional information text
$wsUC = $wb.Worksheets.Item("Universal Credit")
$wsUC.Range("B8").Value = "The demand measure data does not include in-work households that would be eligible for UC. In the analysis, areas with higher than estimated uptake may be influenced by high uptake of UC amongst in-work households.`n`nThere is a lag between the data sourced from the census and the UC data, therefore rates are not exact.`n`nAll data is publicly available."

Write-Output "edit complete"
